$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink currently attached to the URL cell (A56) before any rows shift.
$ws.Range("A56").Hyperlinks.Delete()

# Insert a new blank row above "U.S. CENSUS BUREAU, 2012" (currently row 55).
# This pushes everything below it down by one row:
#   55 (U.S. CENSUS BUREAU, 2012) -> 56
#   56 (the URL, ex-hyperlink)    -> 57
#   57 (blank source cell)        -> 58
#   60 (USCB)                     -> 61
#   61 (long citation sentence)   -> 62
$ws.Rows.Item(55).Insert()

# The newly inserted row 55 becomes a blank line styled like the surrounding
# italic "source" text.
$ws.Range("A55").Font.Italic = $true

# Swap the contents of rows 57 and 58: the blank line moves to 57 and the URL
# text (now bare, no hyperlink) moves to 58.
$ws.Range("A57").Value = ""
$ws.Range("A57").Font.Italic = $true

$ws.Range("A58").Value = "http://www.census.gov/econ/islandareas/"
$ws.Range("A58").Font.Italic = $true
$ws.Range("A58").Font.Underline = $false

# Row 62 (the long citation sentence, shifted down from row 61) is replaced
# with a short "USCB" label, still in the italic "source" style.
$ws.Range("A62").Value = "USCB"
$ws.Range("A62").Font.Italic = $true
